$wb = $excel.ActiveWorkbook

# --- Property sheet edits ---
$ws = $wb.Worksheets.Item("Property")

# Rename header E1 from "newAddress" to "address"
$ws.Range("E1").Value = "address"

# Update existing row 2 values
$ws.Range("B2").Value = "NewPropertyNameByRocky"
$ws.Range("C2").Value = "NewPropertyDescription"
$ws.Range("E2").Value = 238

# Add new rows worth of address data (foreach-style dump down column E)
$ws.Range("E3").Value = "Botany Road"
$ws.Range("E4").Value = 2013
$ws.Range("E5").Value = "Auckland"

# Column width adjustments matching final layout (values chosen so the
# engine's internal character-width rounding lands as close as possible
# to the target widths of 17.375 / 21.125 / 19.75)
$ws.Columns.Item(2).ColumnWidth = 16.714285714285715
$ws.Columns.Item(3).ColumnWidth = 20.428571428571427
$ws.Columns.Item(5).ColumnWidth = 19

# Activate Property tab (becomes the selected sheet; Register's tabSelected
# flag is cleared automatically as a side effect of the activation)
$ws.Activate()

# Move the selection / active cell on the Property sheet
$ws.Range("D5").Select()
